$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

$data = @(
    @("AR205303", "PreIT4"),
    @("AR207371", "F11A"),
    @("R101120", "PreIT4"),
    @("R130071", "PreIT4"),
    @("R570238", "F9B"),
    @("R57087", "PreIT4"),
    @("R572229", "F9A"),
    @("R63317", "PreIT4"),
    @("R64761", "PreIT4"),
    @("R73755", "PreIT4"),
    @("RE596442", "F9B"),
    @("TR117289", "F11A"),
    @("AR206143", "F9A"),
    @("AR207372", "F11A"),
    @("R112072", "PreIT4"),
    @("R268862", "PreIT4"),
    @("R570239", "F9A"),
    @("R57103", "PreIT4"),
    @("R572543", "F9A"),
    @("R63334", "PreIT4"),
    @("R65224", "PreIT4"),
    @("R79139", "PreIT4"),
    @("TR114463", "7RMY20"),
    @("AR207369", "F11A"),
    @("AR207374", "F9A"),
    @("R113154", "PreIT4"),
    @("R44159", "PreIT4"),
    @("R570743", "PreIT4"),
    @("R572071", "F9B"),
    @("R58027", "PreIT4"),
    @("R63400", "PreIT4"),
    @("R65328", "PreiT4"),
    @("R89462", "PreiT4"),
    @("TR117219", "F11A"),
    @("AR206144", "f9B"),
    @("AR207373", "F11A"),
    @("R112910", "PreIT4"),
    @("R41036", "PreIT4"),
    @("R57044", "PreIT4"),
    @("R571830", "F9A"),
    @("R57923", "PreIT4"),
    @("R63338", "PreIT4"),
    @("R65225", "PreIT4"),
    @("R82172", "PreIT4"),
    @("TR116749", "F11A"),
    @("AR207370", "F11A"),
    @("AR207375", "F9B"),
    @("R119598", "PreIT4"),
    @("R569302", "F9A"),
    @("R570744", "PreIT4"),
    @("R572164", "F9A"),
    @("R59494", "PreIT4"),
    @("R63401", "PreIT4"),
    @("R71230", "PreIT4"),
    @("RE590745", "iT4"),
    @("TR117220", "F11A"),
)

$startRow = 2478
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
}

# Resize the query table to include the new rows
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:B2533"))

# Update the ExternalData_1 defined name to reflect the new range
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet2!ExternalData_1") {
        $n.RefersTo = "=Sheet2!`$A`$1:`$B`$2533"
    }
}

# Update scroll/view position to show the newly added rows
$ws.Activate()
$ws.Range("D2378").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 2489
